# 04_figure4.pptx edit: bump the cached "datetimeFigureOut" placeholder
# text on the master + every layout from 2022-01-18 -> 2022-01-19, and
# nudge a handful of axis-label textboxes on slide 1 (figure 4 is "done").

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder text: "2022-01-18" -> "2022-01-19"
#    (slide master + all slide layouts share the same cached field text)
# ---------------------------------------------------------------------
function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePh = $false
        if ($sh.HasTextFrame) {
            try {
                if ($sh.PlaceholderFormat.Type -eq 16) {
                    $isDatePh = $true
                }
            } catch {
                $isDatePh = $false
            }
        }
        if ($isDatePh) {
            if ($sh.TextFrame.TextRange.Text -eq "2022-01-18") {
                $sh.TextFrame.TextRange.Text = "2022-01-19"
            }
        }
    }
}

Update-DatePlaceholder $p.SlideMaster.Shapes

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}

# ---------------------------------------------------------------------
# 2) Slide 1: reposition a few (rotated) axis-label textboxes that live
#    inside the nested "Group 30" shape grouping.
# ---------------------------------------------------------------------
function Get-ItemByName($collection, $name) {
    for ($i = 1; $i -le $collection.Count; $i++) {
        $item = $collection.Item($i)
        if ($item.Name -eq $name) {
            return $item
        }
    }
    return $null
}

$s = $p.Slides.Item(1)
$topGroup = $s.Shapes.Item(1)
$items = $topGroup.GroupItems

$moves = @(
    @{ Name = "TextBox 4";  Left = 80.34;              Top = 262.2996062992126 },
    @{ Name = "TextBox 6";  Left = 384.3492125984252;  Top = 262.59 },
    @{ Name = "TextBox 8";  Left = 719.1464566929134;  Top = 262.3118897637795 },
    @{ Name = "TextBox 9";  Left = 547.892283464567;   Top = 136.80015748031497 },
    @{ Name = "TextBox 12"; Left = -90.68472440944882; Top = 391.85488188976376 },
    @{ Name = "TextBox 19"; Left = 239.19031496062993; Top = 395.0752755905512 },
    @{ Name = "TextBox 21"; Left = 555.4093700787402;  Top = 397.0475590551181 }
)

foreach ($mv in $moves) {
    $shape = Get-ItemByName $items $mv.Name
    if ($shape -ne $null) {
        $shape.Left = $mv.Left
        $shape.Top = $mv.Top
    }
}
